$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 149, shifting existing rows 149-170 down to 150-171
$ws.Rows(149).Insert()

# Copy the date formatting from the (now shifted) next row so the new date cell
# keeps the same "$/caja" style / date number format as the rest of column D
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat

# Populate the new row 149 with the new weekly data point
$ws.Range("A149").Value = 4
$ws.Range("B149").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C149").Value = "Los Lagos"
$ws.Range("D149").Value = 44522
$ws.Range("E149").Value = 10
$ws.Range("F149").Value = 100112021
$ws.Range("G149").Value = "Ají"
$ws.Range("H149").Value = "Inferno"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 40
$ws.Range("K149").Value = 25000
$ws.Range("L149").Value = 25000
$ws.Range("M149").Value = 25000
$ws.Range("N149").Value = "$/caja 12 kilos"
$ws.Range("O149").Value = "Región de Arica y Parinacota"
$ws.Range("P149").Value = 2083
$ws.Range("Q149").Value = 12
$ws.Range("R149").Value = "Hortaliza"
